$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append new row 21 with the new mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(21, 1).Value = "Factuur verzoek"
$logs.Cells.Item(21, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(21, 3).Value = "Kunt u mij de factuur van mijn laatste bestelling toesturen?"
$logs.Cells.Item(21, 4).Value = "Factuur / Administratie"
$logs.Cells.Item(21, 6).Value = "2025-06-19 21:38:18"
$logs.Cells.Item(21, 7).Value = "Nee"

# Extend the conditional formatting ranges to cover the newly added row 21
$catFormats = $logs.Range("D2:D20").FormatConditions
$catFormats.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))

$answeredFormats = $logs.Range("G2:G20").FormatConditions
$answeredFormats.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))

# --- Sheet "Dashboard": refresh the category summary (counts + ordering) ---
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Cells.Item(4, 1).Value = "Factuur / Administratie"
$dashboard.Cells.Item(5, 1).Value = "IT / Technisch probleem"
$dashboard.Cells.Item(7, 1).Value = "Afmelding / Nieuwsbrief"
$dashboard.Cells.Item(7, 2).Value = 2
$dashboard.Cells.Item(8, 1).Value = "Openingstijden / Locatie"
